$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that now hold numeric-looking text must be forced to Text format
# before assignment (otherwise Excel auto-converts "243.51" -> 243.51 as a
# real number and mangles fixed trailing/leading zeros), then the temporary
# Text number-format is cleared again so the cell keeps its original
# (default/General) style - only the stored value changes.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "243.51"
Set-TextValue $ws.Range("D3") "23.71"
Set-TextValue $ws.Range("D4") "3.578"
Set-TextValue $ws.Range("D5") "5.296"
Set-TextValue $ws.Range("D6") "0.05796"
Set-TextValue $ws.Range("D7") "6.481"
Set-TextValue $ws.Range("D8") "3.343"
Set-TextValue $ws.Range("D10") "0.8763"
Set-TextValue $ws.Range("D12") "0.1383"
Set-TextValue $ws.Range("D13") "0.07282"
Set-TextValue $ws.Range("D14") "0.03075"
Set-TextValue $ws.Range("D15") "0.03060"
Set-TextValue $ws.Range("D16") "0.09321"
Set-TextValue $ws.Range("D17") "3.874"
Set-TextValue $ws.Range("D18") "0.001548"
Set-TextValue $ws.Range("D19") "0.04725"
Set-TextValue $ws.Range("D20") "0.006000"
Set-TextValue $ws.Range("D21") "0.001268"
Set-TextValue $ws.Range("D22") "0.004594"
Set-TextValue $ws.Range("D23") "0.00008703"
Set-TextValue $ws.Range("D24") "2.141"
Set-TextValue $ws.Range("D25") "0.3212"
Set-TextValue $ws.Range("D26") "0.1311"
Set-TextValue $ws.Range("D28") "0.0002345"
Set-TextValue $ws.Range("D40") "0.03771"
Set-TextValue $ws.Range("D41") "0.006352"
Set-TextValue $ws.Range("D43") "0.1052"
Set-TextValue $ws.Range("D44") "0.006859"
$ws.Range("E44").Value = "43LocalTradersLCTWorstin24h"
Set-TextValue $ws.Range("D45") "0.00005472"
Set-TextValue $ws.Range("D47") "0.5502"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
Set-TextValue $ws.Range("D48") "0.006420"
Set-TextValue $ws.Range("D49") "0.00002101"
Set-TextValue $ws.Range("D50") "0.0002001"
